# Percentage_of_Instruction.xlsx — "Added seconds of instruction to the
# percent of instruction sheet"
#
# On the VOCALS sheet:
#   - A6 label "Closest" -> "Seconds"
#   - B6/C6/D6 "closest" sample counts replaced with seconds counts
#   - new row 16: "Seconds, (0 dB - 80 dB)" label + three second counts
#     (one real number formatted like text, two numbers stored as text)
#   - columns B:D widened to fit the new numbers
#   - selection left on A15 (just above the new row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VOCALS")

# --- existing "Closest" row -> "Seconds" row ------------------------------
$ws.Range("A6").Value = "Seconds"
$ws.Range("B6").Value = 470
$ws.Range("C6").Value = 1331
$ws.Range("D6").Value = 789

# --- new row 16 -------------------------------------------------------------
# C16: "1060" stored as text, right aligned, numFmt "# ?/?"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "1060"
$ws.Range("C16").NumberFormat = "# ?/?"
$ws.Range("C16").HorizontalAlignment = -4152

# D16: "854" stored as text, right aligned, numFmt "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "854"
$ws.Range("D16").HorizontalAlignment = -4152

# A16: row label (new shared string)
$ws.Range("A16").Value = "Seconds, (0 dB - 80 dB)"

# B16: 455 kept as a real number, right aligned, numFmt "@"
$ws.Range("B16").Value = 455
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").HorizontalAlignment = -4152

# --- column widths so the new numbers fit ----------------------------------
$ws.Range("B:B").ColumnWidth = 9.5
$ws.Range("C:C").ColumnWidth = 10.5
$ws.Range("D:D").ColumnWidth = 9.5

# --- leave the selection where the author left it --------------------------
[void]$ws.Range("A15").Select()
